$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.862.79"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.599.15"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.479"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.246"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.53%  "

$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.821.01"
$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").Value = "1.613.84"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.508"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.40%  "

$ws.Range("D16").Value = "25.855.32"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").Value = "0.0₃0715"
$ws.Range("E18").Value = "  -4.16%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.23%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.129"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.11%  "

$ws.Range("E29").Value = "  -2.09%  "

$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("E31").Value = "  -4.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "

$ws.Range("E33").Value = "  -5.13%  "

$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("D36").Value = "1.104.22"
$ws.Range("E36").Value = "  -2.71%  "

$ws.Range("E37").Value = "  -3.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.796"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0151"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.496"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.31%  "

$ws.Range("D42").Value = "1.733.76"
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("E44").Value = "  -5.09%  "

$ws.Range("D45").Value = "0.0₆0102"
$ws.Range("E45").Value = "  -10.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "53.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.76%  "

$ws.Range("E47").Value = "  -3.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.87%  "
